$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential disclaimer text with the new "as of" date
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-05 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-39
$rows = @(
    @{Row=2; D=0.05878978099067499; E=0.001955416503715224},
    @{Row=3; D=0.05299640687826161; E=-0.00532709148876076},
    @{Row=4; D=0.3055650611686201; E=0.004409603135717699},
    @{Row=5; D=0.03541652408247069; E=-0.01247935456403781},
    @{Row=6; D=0.03225671571563008; E=-0.01215360082747918},
    @{Row=7; D=0.02992813756523227; E=0.01312065860560852},
    @{Row=8; D=0.02870562066900505; E=-0.004172378851999836},
    @{Row=9; D=0.02407733766789289; E=-0.0008527572484366974},
    @{Row=10; D=0.02466881255881598; E=0.003441952809699833},
    @{Row=11; D=0.02383137992984858; E=-0.01049126774720444},
    @{Row=12; D=0.02279921293523826; E=0.009512195121951228},
    @{Row=13; D=0.02167368893831808; E=-0.01487109905020356},
    @{Row=14; D=0.02158118744118402; E=-0.008027352460234782},
    @{Row=15; D=0.02135148430148002; E=-0.002163656579619366},
    @{Row=16; D=0.02199396868851056; E=0.002820051538872992},
    @{Row=17; D=0.02009955941483446; E=-0.01800962996461941},
    @{Row=18; D=0.01460347335101377; E=-0.0008787346221439396},
    @{Row=19; D=0.01695654033706904; E=-0.003884866678439014},
    @{Row=20; D=0.01566643853195312; E=0.001535836177474437},
    @{Row=21; D=0.0164571391906921; E=0.03007264740665661},
    @{Row=22; D=0.01440670716057832; E=-0.003948931116389498},
    @{Row=23; D=0.01505304132089999; E=-0.002585888437384609},
    @{Row=24; D=0.01472388570450851; E=0.0003123048094939573},
    @{Row=25; D=0.01355708358285568; E=0.008282390061131961},
    @{Row=26; D=0.01382646077508769; E=0.008252509783903461},
    @{Row=27; D=0.01266094191119856; E=0.005017103762827935},
    @{Row=28; D=0.01361664813072119; E=0.02685891998869105},
    @{Row=29; D=0.01452540850372145; E=0.0005006257822277593},
    @{Row=30; D=0.01343656429121396; E=0.008404431427479819},
    @{Row=31; D=0.01231585251090769; E=-0.0004862461794944428},
    @{Row=32; D=0.01339828043459663; E=0.01650570676031604},
    @{Row=33; D=0.01218934468303533; E=0.02250296091590998},
    @{Row=34; D=0.006138784327145179; E=0.007473216618761658},
    @{Row=35; D=0.005380913679527762; E=-0.01411025875432259},
    @{Row=36; D=0.005350329369492686; E=-0.01103293891909174},
    @{Row=37; D=0.005247454872101976; E=-0.008171999184838019},
    @{Row=38; D=0.004753828385661733; E=-0.001597156611328643},
    @{Row=39; D=0.9999999999999999; E=0.001283257763709411}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# Restore sheet protection (unprotected above only to allow the edits)
$ws.Protect()
